$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.994.16"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.235.25"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.99%  "
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("D14").Value = "2.579.02"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "2.277.72"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.823"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.78%  "
$ws.Range("D18").Value = "43.889.66"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "0.0₃0956"
$ws.Range("E19").Value = "  -3.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.72%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.95%  "
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("E32").Value = "  -4.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0806"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.55%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.118"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.107"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.98%  "
$ws.Range("E38").Value = "  -11.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0300"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.56%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Value = "1.738.29"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "85.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.188"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.28%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "14.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.41%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.03%  "
